$d = $word.ActiveDocument

# Grab the last paragraph ("Block 2") and add a new paragraph after it
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()

# The newly created paragraph is now the last one; set its text
$new = $d.Paragraphs.Last
$new.Range.Text = "CIS129-6"
